# Transactions.xlsx update
#
# - Column C header "Transaction Details" -> "Transaction"
# - Two new trailing columns added: E "Statement Date", F "Due Date"
#   (column D "Amount" keeps its place; data row 2 only has values in
#   A-D, so E2/F2 stay empty)
# - Column E gets an auto-fitted width; the active selection ends up on E1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Transaction"
$ws.Range("E1").Value = "Statement Date"
$ws.Range("F1").Value = "Due Date"

$ws.Columns.Item(5).AutoFit()

$ws.Range("E1").Select() | Out-Null
